# edit.ps1 - apply "Update countries & provincias Spain" changes
# - bump the "datos actualizados" timestamp in A1
# - swap 5 pairs of adjacent countries whose case counts leapfrogged
#   each other, carrying each country's own stats to its new row
# - refresh the numeric COVID-19 stats for the affected rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 14:22"

# Row 14: Suiza
$ws.Range("B14").Value = 23612
$ws.Range("C14").Value = 332
$ws.Range("E14").Value = 12886
$ws.Range("G14").Value = 31
$ws.Range("H14").Value = 926

# Row 15: Paises Bajos
$ws.Range("B15").Value = 21762
$ws.Range("C15").Value = 1213
$ws.Range("E15").Value = 19116
$ws.Range("G15").Value = 148
$ws.Range("H15").Value = 2396

# Row 19: Austria
$ws.Range("B19").Value = 13105
$ws.Range("C19").Value = 163
$ws.Range("E19").Value = 7570

# Row 23: Suecia
$ws.Range("B23").Value = 9141
$ws.Range("C23").Value = 722
$ws.Range("E23").Value = 8143
$ws.Range("F23").Value = 719
$ws.Range("G23").Value = 106
$ws.Range("H23").Value = 793

# Row 28: Dinamarca
$ws.Range("A28").Value = "Dinamarca"
$ws.Range("B28").Value = 5635
$ws.Range("C28").Value = 233
$ws.Range("D28").Value = 1736
$ws.Range("E28").Value = 3662
$ws.Range("F28").Value = 120
$ws.Range("G28").Value = 19
$ws.Range("H28").Value = 237

# Row 29: Chile
$ws.Range("A29").Value = "Chile"
$ws.Range("B29").Value = 5546
$ws.Range("D29").Value = 1115
$ws.Range("E29").Value = 4383
$ws.Range("F29").Value = 362
$ws.Range("H29").Value = 48

# Row 59: Croacia
$ws.Range("A59").Value = "Croacia"
$ws.Range("B59").Value = 1407
$ws.Range("C59").Value = 64
$ws.Range("D59").Value = 219
$ws.Range("E59").Value = 1168
$ws.Range("F59").Value = 34
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 20

# Row 60: Marruecos
$ws.Range("A60").Value = "Marruecos"
$ws.Range("B60").Value = 1346
$ws.Range("C60").Value = 71
$ws.Range("D60").Value = 103
$ws.Range("E60").Value = 1147
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 96

# Row 64: Moldavia
$ws.Range("D64").Value = 50
$ws.Range("E64").Value = 1096

# Row 69: Lituania
$ws.Range("E69").Value = 931
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 16

# Row 72: Barein
$ws.Range("A72").Value = "Barein"
$ws.Range("B72").Value = 855
$ws.Range("C72").Value = 32
$ws.Range("D72").Value = 495
$ws.Range("E72").Value = 355
$ws.Range("F72").Value = 3
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 5

# Row 73: Bosnia y Herzegovina
$ws.Range("A73").Value = "Bosnia y Herzegovina"
$ws.Range("B73").Value = 841
$ws.Range("C73").Value = 37
$ws.Range("D73").Value = 95
$ws.Range("E73").Value = 711
$ws.Range("F73").Value = 4
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 35

# Row 79: Republica de Macedonia
$ws.Range("A79").Value = "Republica de Macedonia"
$ws.Range("B79").Value = 663
$ws.Range("C79").Value = 46
$ws.Range("D79").Value = 37
$ws.Range("E79").Value = 596
$ws.Range("F79").Value = 15
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 30

# Row 80: Tunez
$ws.Range("A80").Value = "Tunez"
$ws.Range("B80").Value = 628
$ws.Range("D80").Value = 25
$ws.Range("E80").Value = 579
$ws.Range("F80").Value = 67
$ws.Range("H80").Value = 24

# Row 109: Vietnam
$ws.Range("A109").Value = "Vietnam"
$ws.Range("B109").Value = 255
$ws.Range("D109").Value = 128
$ws.Range("E109").Value = 127
$ws.Range("F109").Value = 8
$ws.Range("H109").Value = 0

# Row 110: Montenegro
$ws.Range("A110").Value = "Montenegro"
$ws.Range("B110").Value = 252
$ws.Range("C110").Value = 4
$ws.Range("D110").Value = 4
$ws.Range("E110").Value = 246
$ws.Range("F110").Value = 7
$ws.Range("H110").Value = 2

# Row 111: Senegal
$ws.Range("B111").Value = 250
$ws.Range("C111").Value = 6
$ws.Range("D111").Value = 123
$ws.Range("E111").Value = 125
